$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 43399.555
$ws.Range("J70").Value = 59273.617
$ws.Range("L70").Value = 177820.851
$ws.Range("N70").Value = -178360.851

$ws.Range("H73").Value = 43399.555
$ws.Range("J73").Value = 59273.617
$ws.Range("L73").Value = 177820.851
$ws.Range("N73").Value = -179692.851

$ws.Range("H80").Value = 346.69565
$ws.Range("I80").Value = 319.07693
$ws.Range("J80").Value = 382.6
$ws.Range("K80").Value = 957.2307900000001
$ws.Range("L80").Value = 1147.8
$ws.Range("M80").Value = 40.76920999999993
$ws.Range("N80").Value = -3143.8

$ws.Range("H83").Value = 346.69565
$ws.Range("I83").Value = 319.07693
$ws.Range("J83").Value = 382.6
$ws.Range("K83").Value = 2871.69237
$ws.Range("L83").Value = 3443.4
$ws.Range("M83").Value = 2120.30763
$ws.Range("N83").Value = -13427.4

$ws.Range("H94").Value = 30000
$ws.Range("I94").Value = 30000
$ws.Range("K94").Value = 30000
$ws.Range("M94").Value = -29549

$ws.Range("H108").Value = 1000000000
$ws.Range("I108").Value = 1000000000
$ws.Range("K108").Value = 1000000000
$ws.Range("M108").Value = -999996160

$ws.Range("H137").Value = 2467.4333
$ws.Range("I137").Value = 1261.4445
$ws.Range("J137").Value = 4276.4165
$ws.Range("K137").Value = 3784.3335
$ws.Range("L137").Value = 12829.2495
$ws.Range("M137").Value = -1234.3335
$ws.Range("N137").Value = -17929.2495

$ws.Range("H138").Value = 4079.9788
$ws.Range("J138").Value = 4910.722
$ws.Range("L138").Value = 14732.166
$ws.Range("N138").Value = -25012.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4779.3335
$ws.Range("J61").Value = 4806.5
$ws.Range("L61").Value = 4806.5
$ws.Range("N61").Value = -5230.5

$ws.Range("H136").Value = 4779.3335
$ws.Range("J136").Value = 4806.5
$ws.Range("L136").Value = 14419.5
$ws.Range("N136").Value = -19519.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 811.19354
$ws.Range("I94").Value = 811.19354
$ws.Range("K94").Value = 811.19354
$ws.Range("M94").Value = -360.19354

$ws.Range("H105").Value = 5421.9375
$ws.Range("I105").Value = 3783.4666
$ws.Range("K105").Value = 3783.4666
$ws.Range("M105").Value = -2036.4666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2949.6
$ws.Range("I31").Value = 1557.7142
$ws.Range("J31").Value = 6197.3335
$ws.Range("K31").Value = 1557.7142
$ws.Range("L31").Value = 6197.3335
$ws.Range("M31").Value = -1262.7142
$ws.Range("N31").Value = -6787.3335

$ws.Range("H34").Value = 2949.6
$ws.Range("I34").Value = 1557.7142
$ws.Range("J34").Value = 6197.3335
$ws.Range("K34").Value = 1557.7142
$ws.Range("L34").Value = 6197.3335
$ws.Range("M34").Value = -1355.7142
$ws.Range("N34").Value = -6601.3335

$ws.Range("H107").Value = 2157.5334
$ws.Range("I107").Value = 2341.4348
$ws.Range("J107").Value = 1553.2858
$ws.Range("K107").Value = 2341.4348
$ws.Range("L107").Value = 1553.2858
$ws.Range("M107").Value = -421.4348
$ws.Range("N107").Value = -5393.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2160
$ws.Range("J81").Value = 2990
$ws.Range("L81").Value = 8970
$ws.Range("N81").Value = -11216

$ws.Range("H84").Value = 2160
$ws.Range("J84").Value = 2990
$ws.Range("L84").Value = 26910
$ws.Range("N84").Value = -38142

$ws.Range("H107").Value = 765.75
$ws.Range("I107").Value = 568.6667
$ws.Range("J107").Value = 811.2308
$ws.Range("K107").Value = 1706.0001
$ws.Range("L107").Value = 2433.6924
$ws.Range("M107").Value = 213.9999
$ws.Range("N107").Value = -6273.6924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 58333
$ws.Range("I63").Value = 67999
$ws.Range("J63").Value = 53500
$ws.Range("K63").Value = 67999
$ws.Range("L63").Value = 53500
$ws.Range("M63").Value = -67313
$ws.Range("N63").Value = -54872

$ws.Range("H66").Value = 58333
$ws.Range("I66").Value = 67999
$ws.Range("J66").Value = 53500
$ws.Range("K66").Value = 203997
$ws.Range("L66").Value = 160500
$ws.Range("M66").Value = -200565
$ws.Range("N66").Value = -167364

$ws.Range("H132").Value = 2523.6667
$ws.Range("I132").Value = 1279
$ws.Range("J132").Value = 6257.6665
$ws.Range("K132").Value = 3837
$ws.Range("L132").Value = 18772.9995
$ws.Range("M132").Value = -1307
$ws.Range("N132").Value = -23832.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

$ws.Range("H55").Value = 367.5
$ws.Range("I55").Value = 380.91666
$ws.Range("J55").Value = 287
$ws.Range("K55").Value = 380.91666
$ws.Range("L55").Value = 287
$ws.Range("M55").Value = -207.91666
$ws.Range("N55").Value = -633

$ws.Range("H68").Value = 3748.8333
$ws.Range("J68").Value = 3898.6
$ws.Range("L68").Value = 3898.6
$ws.Range("N68").Value = -5396.6

$ws.Range("H71").Value = 3748.8333
$ws.Range("J71").Value = 3898.6
$ws.Range("L71").Value = 19493
$ws.Range("N71").Value = -26981

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H136").Value = 500
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 68000
$ws.Range("J64").Value = 68000
$ws.Range("L64").Value = 68000
$ws.Range("N64").Value = -68496

$ws.Range("H67").Value = 68000
$ws.Range("J67").Value = 68000
$ws.Range("L67").Value = 68000
$ws.Range("N67").Value = -69716

$ws.Range("H107").Value = 478.95
$ws.Range("I107").Value = 387.83334
$ws.Range("K107").Value = 1163.50002
$ws.Range("M107").Value = 756.4999800000001

$ws.Range("H122").Value = 6942
$ws.Range("I122").Value = 7330.4
$ws.Range("K122").Value = 21991.2
$ws.Range("M122").Value = -19541.2

$ws.Range("H126").Value = 1609.1111
$ws.Range("I126").Value = 1110.6
$ws.Range("K126").Value = 3331.8
$ws.Range("M126").Value = -861.7999999999997
